$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded ahead of the existing
# history for "Terminal La Palmera de La Serena - Berenjena": insert a
# fresh row at row 140 (pushing the former rows 140-147 down to 141-148)
# and populate it with the new reading.
$ws.Rows(140).Insert()

$ws.Cells.Item(140, 1).Value = 8
$ws.Cells.Item(140, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(140, 3).Value = "Coquimbo"
$ws.Cells.Item(140, 4).Value = 44753
$ws.Cells.Item(140, 5).Value = 4
$ws.Cells.Item(140, 6).Value = 100112001
$ws.Cells.Item(140, 7).Value = "Berenjena"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 460
$ws.Cells.Item(140, 11).Value = 10500
$ws.Cells.Item(140, 12).Value = 11000
$ws.Cells.Item(140, 13).Value = 10750
$ws.Cells.Item(140, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(140, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(140, 16).Value = 215
$ws.Cells.Item(140, 17).Value = 50
$ws.Cells.Item(140, 18).Value = "Hortaliza"
